$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 2 (pushes existing data rows 2-22 down to 3-23)
$ws.Rows.Item(2).Insert()

# Populate the new row 2 with the newly collected accelerometer sample
$ws.Range("A2").Value = -3.729709470272064
$ws.Range("B2").Value = 9.457800364494323
$ws.Range("C2").Value = 0.187229474633932

# The insert pushed the original last two rows (old rows 21 and 22, now at
# rows 22 and 23) past the new end of the dataset; remove them so the sheet
# again ends at row 21 (dimension A1:C21).
$ws.Rows.Item(22).Delete()
$ws.Rows.Item(22).Delete()
